# All demo tables are now exported to Sexy. Add a "Sexy Header" row to the
# "#Users" table so tables don't have to script to Sexy by hand.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#Users")

# Write the value column first, then the label column, so the shared
# strings table picks up the two new strings in the same order as Excel
# would (path before the friendly name).
$ws.Range("B13").Value = "tables\rococo.tables.test.sxh"
$ws.Range("A13").Value = "Sexy Header"

# Move the active selection down to the next empty row, matching what
# Excel leaves behind after typing a new row of data.
$ws.Activate() | Out-Null
$ws.Range("A14").Select() | Out-Null
